# Insert a new weekly data row for Pomelo (Start Ruby, Primera) at row 8,
# shifting the existing rows 8-38 down to rows 9-39.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = Get-Date -Year 2022 -Month 7 -Day 7 -Hour 0 -Minute 0 -Second 0
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100102
$ws.Range("H8").Value = "Cítricos"
$ws.Range("I8").Value = 100102006
$ws.Range("J8").Value = "Pomelo"
$ws.Range("K8").Value = "Start Ruby"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 140000
$ws.Range("O8").Value = 140000
$ws.Range("P8").Value = 140000
$ws.Range("Q8").Value = '$/bins (350 kilos)'
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 400
$ws.Range("T8").Value = 350
